# "6 february Presenti sheet" - mark Feb 6 attendance on the Feb-2024 sheet
# and refresh the reminder textbox / selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feb-2024")
$ws.Activate()

# Columns F, G, H on row 1 are 4-Feb-2024, 5-Feb-2024 and 6-Feb-2024.
# Row 2 = Pratiksha Bhuse(TL), Row 3 = Sangita Survase.
# 4th & 5th were Absent, 6th (today) is Present for both students.
$ws.Range("F2").Value = "Absent"
$ws.Range("G2").Value = "Absent"
$ws.Range("H2").Value = "Present"

$ws.Range("F3").Value = "Absent"
$ws.Range("G3").Value = "Absent"
$ws.Range("H3").Value = "Present"

# Update the reminder textbox: drop the old "Apptitude" line and bump the
# class time from 6 p.m to 7 p.m.
$shp = $ws.Shapes.Item("TextBox 1")
$shp.TextFrame.Characters().Text = "  Time: 5 p.m to 7 p.m"

# Move the selection cursor to where the user last clicked.
$ws.Range("I12").Select()
